$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that must stay as literal text
# (matching the source workbook's inlineStr cells, e.g. "1.00", "0.0000191").
# Force Text format first so Excel does not coerce them to numbers, then
# clear the formatting delta afterwards so no stray style is introduced.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "63.756.32"
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("D3").Value = "3.430.34"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "572.66"
$ws.Range("E5").Value = "  +2.98%  "
$ws.Range("D6").Value = "156.68"
$ws.Range("E6").Value = "  +2.72%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.432.01"
$ws.Range("E8").Value = "  +2.43%  "
$ws.Range("E9").Value = "  +2.77%  "
$ws.Range("D10").Value = "7.45"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("E11").Value = "  +3.61%  "
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "4.023.23"
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("E14").Value = "  -3.14%  "
$ws.Range("D15").Value = "0.0000191"
$ws.Range("E15").Value = "  +5.20%  "
$ws.Range("D16").Value = "27.20"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").Value = "63.911.20"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").Value = "3.381.72"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("E19").Value = "  -2.25%  "
$ws.Range("D20").Value = "14.22"
$ws.Range("E20").Value = "  +3.31%  "
$ws.Range("D21").Value = "388.07"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "8.26"
$ws.Range("E22").Value = "  -2.20%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "72.33"
$ws.Range("E23").Value = "  +2.34%  "
$ws.Range("D24").Value = "0.538"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "0.995"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").Value = "0.0000120"
$ws.Range("E26").Value = "  +23.29%  "
$ws.Range("D27").Value = "9.56"
$ws.Range("E27").Value = "  +8.47%  "
$ws.Range("D28").Value = "0.177"
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").Value = "6.05"
$ws.Range("E30").Value = "  +8.23%  "
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("D33").Value = "6.47"
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("D34").Value = "23.35"
$ws.Range("E34").Value = "  +1.35%  "
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "6.89"
$ws.Range("E36").Value = "  +2.69%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "159.26"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "1.47"
$ws.Range("E38").Value = "  -0.82%  "
$ws.Range("D39").Value = "0.0769"
$ws.Range("E39").Value = "  +2.88%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.84"
$ws.Range("E40").Value = "  -2.35%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.910.27"
$ws.Range("E41").Value = "  +3.08%  "
$ws.Range("D42").Value = "26.98"
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("D43").Value = "0.0319"
$ws.Range("E43").Value = "  +2.54%  "
$ws.Range("D44").Value = "4.41"
$ws.Range("E44").Value = "  +1.87%  "
$ws.Range("D45").Value = "0.766"
$ws.Range("E45").Value = "  +2.77%  "
$ws.Range("D46").Value = "41.40"
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("D47").Value = "23.70"
$ws.Range("E47").Value = "  +7.66%  "
$ws.Range("D48").Value = "1.08"
$ws.Range("E48").Value = "  +4.14%  "
$ws.Range("E49").Value = "  +22.10%  "
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").Value = "0.848"
$ws.Range("E50").Value = "  +5.12%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "6.47"
$ws.Range("E51").Value = "  +2.77%  "

$dRange.ClearFormats()
